# "adding the employer part" - append newly-classified resume files to the
# DataBase sheet (Path / Class columns), fixing up the previous row's
# classification and appending five new rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (my_files\lspan1np.docx) gets reclassified.
$ws.Range("B10").Value = "GraphicsAndDesign"

$ws.Range("A11").Value = "my_files\f0dybxje.pdf"
$ws.Range("B11").Value = "בס -- 0 --> ד -- Attribute not found"

$ws.Range("A12").Value = "my_files\awqiqrmr.pdf"
$ws.Range("B12").Value = "בס -- 0 --> ד -- Attribute not found"

$ws.Range("A13").Value = "my_files\2enid42m.pdf"
$ws.Range("B13").Value = "בס -- 0 --> ד -- 0 --> קורות -- "

$ws.Range("A14").Value = "my_files\idpfilow.doc"
$ws.Range("B14").Value = "בס -- 0 --> ד -- "

$ws.Range("A15").Value = "my_files\kzyqscg2.doc"
$ws.Range("B15").Value = "בס -- Attribute not found"

# Leave the selection on the row that was just edited, as in the workbook.
$ws.Range("A10").Select()
